$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.975.83"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "2.512.27"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("E5").Value = "  -1.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.66%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.74%  "

$ws.Range("D9").Value = "2.513.31"
$ws.Range("E9").Value = "  -0.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1000"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.64%  "

$ws.Range("E12").Value = "  -2.89%  "

$ws.Range("E13").Value = "  -0.16%  "

$ws.Range("D14").Value = "2.955.98"
$ws.Range("E14").Value = "  +0.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.58%  "

$ws.Range("D16").Value = "58.899.00"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").Value = "2.511.00"
$ws.Range("E18").Value = "  -0.40%  "

$ws.Range("E19").Value = "  -1.97%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "321.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.31%  "

$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.422"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.35%  "

$ws.Range("E26").Value = "  +1.94%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.52%  "

$ws.Range("E30").Value = "  -1.12%  "

$ws.Range("E31").Value = "  -1.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "161.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.44%  "

$ws.Range("E33").Value = "  +0.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.40"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.56%  "

$ws.Range("E38").Value = "  -1.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.97"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.12%  "

$ws.Range("E40").Value = "  -1.89%  "

$ws.Range("E41").Value = "  -2.53%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "278.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("E45").Value = "  +0.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.594"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.17%  "

$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.73%  "

$ws.Range("E50").Value = "  -0.94%  "

$ws.Range("E51").Value = "  -2.90%  "
